$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-08-26 Monday" "2024-08-27 Tuesday"

Replace-Text "20×41=820" "45×43=1935"
Replace-Text "46×80=3680" "52×76=3952"
Replace-Text "25×71=1775" "93×43=3999"
Replace-Text "48×30=1440" "72×19=1368"
Replace-Text "58×99=5742" "71×68=4828"
Replace-Text "71×32=2272" "65×96=6240"
Replace-Text "30×24=720" "39×12=468"
Replace-Text "94×28=2632" "43×97=4171"
Replace-Text "94×70=6580" "73×23=1679"
Replace-Text "95×61=5795" "11×70=770"
Replace-Text "79×41=3239" "70×43=3010"
Replace-Text "45×83=3735" "45×20=900"
Replace-Text "83×22=1826" "45×18=810"
Replace-Text "32×81=2592" "90×62=5580"
Replace-Text "75×21=1575" "17×87=1479"
Replace-Text "52×82=4264" "89×46=4094"
Replace-Text "94×63=5922" "26×66=1716"
Replace-Text "81×89=7209" "22×22=484"
Replace-Text "16×35=560" "37×50=1850"
Replace-Text "75×54=4050" "90×64=5760"
Replace-Text "28×55=1540" "17×11=187"
Replace-Text "13×88=1144" "95×24=2280"
Replace-Text "16×80=1280" "32×27=864"
Replace-Text "79×91=7189" "19×42=798"
Replace-Text "61×37=2257" "26×72=1872"
